# Add the new code-table worksheets (HousingStatus, IncomeLevel, Occupation,
# Education, Language) at the end of the workbook, matching the author's
# commit "Add spread sheets for the new code tables."
#
# The cell-write order below reproduces the shared-string allocation order
# of the original commit: HousingStatus is filled immediately; the other
# four sheets get their headers typed first, then their bodies are filled
# in the order Education, Occupation, Language, and IncomeLevel last.

$wb = $excel.ActiveWorkbook

function Add-SheetAtEnd($wb, $name) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $ws.Name = $name
    return $ws
}

function Fill-Body($ws, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $i + 1
        $ws.Cells.Item($row, 2).Value = $values[$i]
    }
}

# --- HousingStatus (sheet filled immediately, start to finish) --------
$wsHousing = Add-SheetAtEnd $wb "HousingStatus"
$wsHousing.Cells.Item(1, 1).Value = "HousingStatusID"
$wsHousing.Cells.Item(1, 2).Value = "HousingStatusDescription"
Fill-Body $wsHousing @(
    "Housing Status 1",
    "Housing Status 2",
    "Housing Status 3",
    "Housing Status 4",
    "Housing Status 5",
    "Housing Status 6",
    "Housing Status 7"
)
$wsHousing.Range("C8").Select() | Out-Null

# --- IncomeLevel (header only for now) --------------------------------
$wsIncome = Add-SheetAtEnd $wb "IncomeLevel"
$wsIncome.Cells.Item(1, 1).Value = "IncomeLevelID"
$wsIncome.Cells.Item(1, 2).Value = "IncomeLevel"

# --- Occupation (header only for now) -------------------------------
$wsOccupation = Add-SheetAtEnd $wb "Occupation"
$wsOccupation.Cells.Item(1, 1).Value = "OccupationID"
$wsOccupation.Cells.Item(1, 2).Value = "Occupation"

# --- Education (header + body filled now) ---------------------------
$wsEducation = Add-SheetAtEnd $wb "Education"
$wsEducation.Cells.Item(1, 1).Value = "EducationID"
$wsEducation.Cells.Item(1, 2).Value = "EducationLevel"
Fill-Body $wsEducation @(
    "8th grade or less",
    "Some high school",
    "GED",
    "High school diploma ",
    "Some college",
    "College graduate or more "
)
$wsEducation.Range("C14").Select() | Out-Null

# --- back to Occupation: fill body now --------------------------------
Fill-Body $wsOccupation @(
    "Occupation 1",
    "Occupation 2",
    "Occupation 3",
    "Occupation 4",
    "Occupation 5",
    "Occupation 6",
    "Occupation 7",
    "Occupation 8",
    "Occupation 9",
    "Occupation 10",
    "Occupation 11",
    "Occupation 12"
)
$wsOccupation.Range("D12").Select() | Out-Null

# --- Language (header + body filled now) ------------------------------
$wsLanguage = Add-SheetAtEnd $wb "Language"
$wsLanguage.Cells.Item(1, 1).Value = "LanguageID"
$wsLanguage.Cells.Item(1, 2).Value = "Language"
Fill-Body $wsLanguage @(
    "English",
    "Spanish",
    "Chinese",
    "Tagalog",
    "French",
    "Vietnamese",
    "German",
    "Korean",
    "Russian",
    "Arabic",
    "Italian",
    "Portuguese",
    "Hungarian",
    "Polish ",
    "Hindi ",
    "ASL",
    "Japanese ",
    "Persian",
    "Urdu",
    "Gujarati",
    "Greek",
    "Serbo-Croatian",
    "Punjabi ",
    "Armenian ",
    "Hebrew ",
    "Cambodian",
    "Hmong",
    "Navajo ",
    "Thai",
    "Yiddish ",
    "Laotian "
)
$wsLanguage.Range("N33").Select() | Out-Null

# --- back to IncomeLevel: fill body last -------------------------------
Fill-Body $wsIncome @(
    "None",
    "Less than `$300 ",
    "`$300-599 ",
    "`$600-999 ",
    "`$1,000-1,999 ",
    "`$2,000 or more "
)
$wsIncome.Range("B20").Select() | Out-Null

# Incidental selection change on PersonAge left over from the same editing
# session (D34 -> A34).
$wb.Worksheets.Item("PersonAge").Range("A34").Select() | Out-Null

# The commit left IncomeLevel as the active tab/cell.
$wsIncome.Select()
$wsIncome.Range("B20").Select() | Out-Null
